# Update cryptos list with latest prices, percentage changes, and reordered rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.057.47"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "3.482.05"
$ws.Range("E3").Value = "  -4.10%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'197.97"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").Value = "'548.28"
$ws.Range("D7").Value = "3.472.99"
$ws.Range("E7").Value = "  -4.21%  "
$ws.Range("E8").Value = "  -2.98%  "
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  -4.47%  "
$ws.Range("D11").Value = "'61.96"
$ws.Range("E11").Value = "  +10.98%  "
$ws.Range("E12").Value = "  -7.44%  "
$ws.Range("E13").Value = "  -9.83%  "
$ws.Range("D14").Value = "'9.75"
$ws.Range("D15").Value = "4.045.40"
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("D16").Value = "3.487.72"
$ws.Range("E16").Value = "  -3.71%  "
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").Value = "66.800.90"
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("D19").Value = "'18.19"
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("E20").Value = "  -6.60%  "
$ws.Range("E21").Value = "  -5.96%  "
$ws.Range("D22").Value = "'387.17"
$ws.Range("E22").Value = "  -4.20%  "
$ws.Range("D23").Value = "'3.97"
$ws.Range("E23").Value = "  -6.12%  "
$ws.Range("D24").Value = "'11.79"
$ws.Range("E24").Value = "  -7.14%  "
$ws.Range("D25").Value = "'82.04"
$ws.Range("E25").Value = "  -4.80%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'3.82"
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'12.13"
$ws.Range("E27").Value = "  -3.99%  "
$ws.Range("E28").Value = "  -6.08%  "
$ws.Range("D29").Value = "'8.75"
$ws.Range("E29").Value = "  -4.60%  "
$ws.Range("D30").Value = "'30.83"
$ws.Range("E30").Value = "  -2.94%  "
$ws.Range("D31").Value = "'673.77"
$ws.Range("E31").Value = "  -2.60%  "
$ws.Range("E32").Value = "  -14.74%  "
$ws.Range("D33").Value = "'11.63"
$ws.Range("E33").Value = "  -4.88%  "
$ws.Range("D34").Value = "'63.26"
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("E35").Value = "  -7.60%  "
$ws.Range("D36").Value = "'38.08"
$ws.Range("E36").Value = "  -10.91%  "
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "'0.396"
$ws.Range("E38").Value = "  -5.30%  "
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'0.998"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.062.37"
$ws.Range("E40").Value = "  -2.67%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.129"
$ws.Range("E41").Value = "  -4.94%  "
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").Value = "'2.96"
$ws.Range("E42").Value = "  -5.89%  "
$ws.Range("D43").Value = "0.0₃0669"
$ws.Range("E43").Value = "  -16.48%  "
$ws.Range("D44").Value = "'2.77"
$ws.Range("E44").Value = "  +6.11%  "
$ws.Range("D45").Value = "'2.48"
$ws.Range("E45").Value = "  -13.30%  "
$ws.Range("D46").Value = "'2.72"
$ws.Range("E46").Value = "  -7.60%  "
$ws.Range("E47").Value = "  -7.56%  "
$ws.Range("E48").Value = "  -5.30%  "
$ws.Range("D49").Value = "'136.37"
$ws.Range("E49").Value = "  -4.30%  "
$ws.Range("D50").Value = "'2.94"
$ws.Range("E50").Value = "  -5.48%  "
$ws.Range("D51").Value = "'8.14"
$ws.Range("E51").Value = "  -8.21%  "
